# Insert 3 new weekly price rows at the top of the "Pepino dulce" data block
# (Terminal La Palmera de La Serena), shifting the existing rows 274-286
# down to 277-289, and fill the new rows 274-276 with the latest week's
# data (Fecha 2022-02-18, serial 44610) for qualities Primera/Segunda/Tercera.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 274:286 down by three rows, creating blank rows 274:276.
$ws.Range("A274:R276").Insert()

# Common / repeated field values for the three new rows.
$mercadoId = 8
$mercado   = "Terminal La Palmera de La Serena"
$region    = "Coquimbo"
$fecha     = 44610
$codreg    = 4
$catId     = 100112043
$categoria = "Pepino dulce"
$variedad  = "Cultivar IV Región"
$unidad    = "`$/bandeja 18 kilos"
$origen    = "Provincia de Limarí"
$kgUnidad  = 18
$clasif    = "Hortaliza"

# Row 274: Primera
$ws.Cells.Item(274, 1).Value  = $mercadoId
$ws.Cells.Item(274, 2).Value  = $mercado
$ws.Cells.Item(274, 3).Value  = $region
$ws.Cells.Item(274, 4).Value  = $fecha
$ws.Cells.Item(274, 5).Value  = $codreg
$ws.Cells.Item(274, 6).Value  = $catId
$ws.Cells.Item(274, 7).Value  = $categoria
$ws.Cells.Item(274, 8).Value  = $variedad
$ws.Cells.Item(274, 9).Value  = "Primera"
$ws.Cells.Item(274, 10).Value = 400
$ws.Cells.Item(274, 11).Value = 11500
$ws.Cells.Item(274, 12).Value = 12000
$ws.Cells.Item(274, 13).Value = 11750
$ws.Cells.Item(274, 14).Value = $unidad
$ws.Cells.Item(274, 15).Value = $origen
$ws.Cells.Item(274, 16).Value = 653
$ws.Cells.Item(274, 17).Value = $kgUnidad
$ws.Cells.Item(274, 18).Value = $clasif

# Row 275: Segunda
$ws.Cells.Item(275, 1).Value  = $mercadoId
$ws.Cells.Item(275, 2).Value  = $mercado
$ws.Cells.Item(275, 3).Value  = $region
$ws.Cells.Item(275, 4).Value  = $fecha
$ws.Cells.Item(275, 5).Value  = $codreg
$ws.Cells.Item(275, 6).Value  = $catId
$ws.Cells.Item(275, 7).Value  = $categoria
$ws.Cells.Item(275, 8).Value  = $variedad
$ws.Cells.Item(275, 9).Value  = "Segunda"
$ws.Cells.Item(275, 10).Value = 300
$ws.Cells.Item(275, 11).Value = 9500
$ws.Cells.Item(275, 12).Value = 10000
$ws.Cells.Item(275, 13).Value = 9750
$ws.Cells.Item(275, 14).Value = $unidad
$ws.Cells.Item(275, 15).Value = $origen
$ws.Cells.Item(275, 16).Value = 542
$ws.Cells.Item(275, 17).Value = $kgUnidad
$ws.Cells.Item(275, 18).Value = $clasif

# Row 276: Tercera
$ws.Cells.Item(276, 1).Value  = $mercadoId
$ws.Cells.Item(276, 2).Value  = $mercado
$ws.Cells.Item(276, 3).Value  = $region
$ws.Cells.Item(276, 4).Value  = $fecha
$ws.Cells.Item(276, 5).Value  = $codreg
$ws.Cells.Item(276, 6).Value  = $catId
$ws.Cells.Item(276, 7).Value  = $categoria
$ws.Cells.Item(276, 8).Value  = $variedad
$ws.Cells.Item(276, 9).Value  = "Tercera"
$ws.Cells.Item(276, 10).Value = 200
$ws.Cells.Item(276, 11).Value = 7500
$ws.Cells.Item(276, 12).Value = 8000
$ws.Cells.Item(276, 13).Value = 7750
$ws.Cells.Item(276, 14).Value = $unidad
$ws.Cells.Item(276, 15).Value = $origen
$ws.Cells.Item(276, 16).Value = 431
$ws.Cells.Item(276, 17).Value = $kgUnidad
$ws.Cells.Item(276, 18).Value = $clasif
